$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 30: "Ct" message for end-of-turn notification.
# This pushes the existing Cm (move) request/broadcast rows from 30/31 down
# to 31/32, and the trailing blank rows down by one (44 -> 45), which is
# exactly what the target sheet needs (dimension grows from E44 to E45).
$ws.Rows("30").Insert()

# Copy the formatting (borders / wrap text styles) from the row above so the
# new row 30 matches the rest of the table instead of being unstyled.
$ws.Range("A29:E29").Copy()
$ws.Range("A30:E30").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A30").Value = "Client"
$ws.Range("B30").Value = "Serveur"
$ws.Range("C30").Value = "Ct"
$ws.Range("D30").Value = ""
$ws.Range("E30").Value = "Le client indique qu'il a terminé son tour de jeu."

# Row 31 (old row 30 - "Cm" move request) and row 32 (old row 31 - "Cm" move
# broadcast) already hold the right content after the shift, nothing to do.

# --- Fill the next blank row (33) with the new "Cp" PM-sync message.
$ws.Range("A32:E32").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A33").Value = "Serveur"
$ws.Range("B33").Value = "Client"
$ws.Range("C33").Value = "Cp"
$ws.Range("D33").Value = "idPerso;nbPM"
$ws.Range("E33").Value = "Le serveur envoi l'information du nombre de PM restant du personnage actif (synchronisation quand il y a connexion en cours de combat)"
$ws.Rows("33").RowHeight = 45

# --- Match the author's final selection / scroll position.
$ws.Range("E34").Select()
